{"js": "// Office.js (Word JavaScript API)\n// Replicates: place cursor in the (single, empty) paragraph, type\n// \"This is a test\", then press Enter \u2014 leaving a new empty paragraph\n// after it, matching the target OOXML in the diff.\n\nconst body = context.document.body;\n\n// There is a single, empty paragraph in the document to start with.\nconst firstParagraph = body.paragraphs.getFirst();\n\n// Insert the text into that (empty) paragraph.\nfirstParagraph.insertText(\"This is a test\", Word.InsertLocation.start);\nawait context.sync();\n\n// Press Enter: insert a new empty paragraph right after it.\nfirstParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style)\n# Replicates: place cursor in the (single, empty) paragraph, type\n# \"This is a test\", then press Enter \u2014 leaving a new empty paragraph\n# after it, matching the target OOXML in the diff.\n\n$d = $word.ActiveDocument\n\n# There is a single, empty paragraph in the document to start with.\n$para = $d.Paragraphs.First\n\n# Insert the text into that (empty) paragraph.\n$para.Range.InsertAfter(\"This is a test\")\n\n# Press Enter: insert a new empty paragraph right after it.\n$para.Range.InsertParagraphAfter()\n"}
